$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 612
$ws.Cells.Item(2, 9).Value = 612
$ws.Cells.Item(2, 11).Value = 612
$ws.Cells.Item(2, 13).Value = -499

$ws.Cells.Item(11, 8).Value = 16.625
$ws.Cells.Item(11, 9).Value = 16.625
$ws.Cells.Item(11, 11).Value = 16.625
$ws.Cells.Item(11, 13).Value = 123.375

$ws.Cells.Item(32, 8).Value = 9499
$ws.Cells.Item(32, 10).Value = 10249
$ws.Cells.Item(32, 12).Value = 10249
$ws.Cells.Item(32, 14).Value = -10901

$ws.Cells.Item(132, 8).Value = 3368.8147
$ws.Cells.Item(132, 9).Value = 2693.238
$ws.Cells.Item(132, 11).Value = 8079.714
$ws.Cells.Item(132, 13).Value = -5549.714

$ws.Cells.Item(138, 8).Value = 3102.65
$ws.Cells.Item(138, 9).Value = 1758.375
$ws.Cells.Item(138, 10).Value = 3998.8333
$ws.Cells.Item(138, 11).Value = 5275.125
$ws.Cells.Item(138, 12).Value = 11996.4999
$ws.Cells.Item(138, 13).Value = -135.125
$ws.Cells.Item(138, 14).Value = -22276.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 2700.2222
$ws.Cells.Item(132, 9).Value = 2328.8572
$ws.Cells.Item(132, 11).Value = 6986.571599999999
$ws.Cells.Item(132, 13).Value = -4456.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1499
$ws.Cells.Item(86, 9).Value = 1499
$ws.Cells.Item(86, 11).Value = 1499
$ws.Cells.Item(86, 13).Value = -376

$ws.Cells.Item(89, 8).Value = 1499
$ws.Cells.Item(89, 9).Value = 1499
$ws.Cells.Item(89, 11).Value = 7495
$ws.Cells.Item(89, 13).Value = -1879

$ws.Cells.Item(134, 8).Value = 4280.1816
$ws.Cells.Item(134, 9).Value = 2697.75
$ws.Cells.Item(134, 11).Value = 8093.25
$ws.Cells.Item(134, 13).Value = -5558.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1549.6666
$ws.Cells.Item(16, 9).Value = 1359.6
$ws.Cells.Item(16, 11).Value = 1359.6
$ws.Cells.Item(16, 13).Value = -1072.6

$ws.Cells.Item(22, 8).Value = 499
$ws.Cells.Item(22, 9).Value = 499
$ws.Cells.Item(22, 11).Value = 499
$ws.Cells.Item(22, 13).Value = -149

$ws.Cells.Item(58, 8).Value = 4995
$ws.Cells.Item(58, 9).Value = 4995
$ws.Cells.Item(58, 11).Value = 4995
$ws.Cells.Item(58, 13).Value = -4792

$ws.Cells.Item(105, 8).Value = 2196.6
$ws.Cells.Item(105, 9).Value = 2121
$ws.Cells.Item(105, 10).Value = 2499
$ws.Cells.Item(105, 11).Value = 2121
$ws.Cells.Item(105, 12).Value = 2499
$ws.Cells.Item(105, 13).Value = -374
$ws.Cells.Item(105, 14).Value = -5993

$ws.Cells.Item(113, 8).Value = 1549.6666
$ws.Cells.Item(113, 9).Value = 1359.6
$ws.Cells.Item(113, 11).Value = 1359.6
$ws.Cells.Item(113, 13).Value = 810.4000000000001

$ws.Cells.Item(134, 8).Value = 1056.5
$ws.Cells.Item(134, 9).Value = 1050.2
$ws.Cells.Item(134, 11).Value = 3150.6
$ws.Cells.Item(134, 13).Value = -615.6000000000004

$ws.Cells.Item(136, 8).Value = 4995
$ws.Cells.Item(136, 9).Value = 4995
$ws.Cells.Item(136, 11).Value = 14985
$ws.Cells.Item(136, 13).Value = -12435

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(50, 8).Value = 2643.75
$ws.Cells.Item(50, 9).Value = 537.5
$ws.Cells.Item(50, 11).Value = 1612.5
$ws.Cells.Item(50, 13).Value = -1131.5

$ws.Cells.Item(53, 8).Value = 2643.75
$ws.Cells.Item(53, 9).Value = 537.5
$ws.Cells.Item(53, 11).Value = 1612.5
$ws.Cells.Item(53, 13).Value = -1131.5

$ws.Cells.Item(131, 8).Value = 2516.6667
$ws.Cells.Item(131, 10).Value = 2516.6667
$ws.Cells.Item(131, 12).Value = 7550.000100000001
$ws.Cells.Item(131, 14).Value = -17630.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(58, 8).Value = 21666.334
$ws.Cells.Item(58, 10).Value = 19999.5
$ws.Cells.Item(58, 12).Value = 19999.5
$ws.Cells.Item(58, 14).Value = -20553.5

$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).ClearContents()
$ws.Cells.Item(70, 14).ClearContents()

$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).ClearContents()
$ws.Cells.Item(73, 14).ClearContents()

$ws.Cells.Item(97, 8).Value = 1620
$ws.Cells.Item(97, 10).Value = 145
$ws.Cells.Item(97, 12).Value = 145
$ws.Cells.Item(97, 14).Value = -1137

$ws.Cells.Item(102, 8).Value = 2328.4375
$ws.Cells.Item(102, 9).Value = 2328.4375
$ws.Cells.Item(102, 11).Value = 2328.4375
$ws.Cells.Item(102, 13).Value = -706.4375

$ws.Cells.Item(113, 8).Value = 1666
$ws.Cells.Item(113, 9).Value = 1499.5
$ws.Cells.Item(113, 11).Value = 1499.5
$ws.Cells.Item(113, 13).Value = 670.5

$ws.Cells.Item(132, 8).Value = 7045
$ws.Cells.Item(132, 9).Value = 6249.5
$ws.Cells.Item(132, 11).Value = 18748.5
$ws.Cells.Item(132, 13).Value = -16218.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1620.9524
$ws.Cells.Item(22, 9).Value = 1587.6471
$ws.Cells.Item(22, 10).Value = 1762.5
$ws.Cells.Item(22, 11).Value = 1587.6471
$ws.Cells.Item(22, 12).Value = 1762.5
$ws.Cells.Item(22, 13).Value = -1292.6471
$ws.Cells.Item(22, 14).Value = -2352.5

$ws.Cells.Item(27, 8).Value = 1620.9524
$ws.Cells.Item(27, 9).Value = 1587.6471
$ws.Cells.Item(27, 10).Value = 1762.5
$ws.Cells.Item(27, 11).Value = 1587.6471
$ws.Cells.Item(27, 12).Value = 1762.5
$ws.Cells.Item(27, 13).Value = -1480.6471
$ws.Cells.Item(27, 14).Value = -1976.5

$ws.Cells.Item(46, 8).Value = 4706.121
$ws.Cells.Item(46, 10).Value = 4850.467
$ws.Cells.Item(46, 12).Value = 4850.467
$ws.Cells.Item(46, 14).Value = -5226.467

$ws.Cells.Item(57, 8).Value = 5000
$ws.Cells.Item(57, 9).Value = 5000
$ws.Cells.Item(57, 11).Value = 5000
$ws.Cells.Item(57, 13).Value = -4434

$ws.Cells.Item(61, 8).Value = 7666.6665
$ws.Cells.Item(61, 9).Value = 7666.6665
$ws.Cells.Item(61, 11).Value = 7666.6665
$ws.Cells.Item(61, 13).Value = -7464.6665

$ws.Cells.Item(76, 8).Value = 20684.625
$ws.Cells.Item(76, 10).Value = 20884.428
$ws.Cells.Item(76, 12).Value = 20884.428
$ws.Cells.Item(76, 14).Value = -21560.428

$ws.Cells.Item(79, 8).Value = 20684.625
$ws.Cells.Item(79, 10).Value = 20884.428
$ws.Cells.Item(79, 12).Value = 20884.428
$ws.Cells.Item(79, 14).Value = -23224.428

$ws.Cells.Item(113, 8).Value = 7666.6665
$ws.Cells.Item(113, 9).Value = 7666.6665
$ws.Cells.Item(113, 11).Value = 7666.6665
$ws.Cells.Item(113, 13).Value = -5496.6665

$ws.Cells.Item(125, 8).Value = 38000
$ws.Cells.Item(125, 10).Value = 38000
$ws.Cells.Item(125, 12).Value = 38000
$ws.Cells.Item(125, 14).Value = -47840

$ws.Cells.Item(132, 8).Value = 6498.5
$ws.Cells.Item(132, 9).Value = 6498.5
$ws.Cells.Item(132, 11).Value = 19495.5
$ws.Cells.Item(132, 13).Value = -16965.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(93, 8).Value = 4000
$ws.Cells.Item(93, 9).Value = 4000
$ws.Cells.Item(93, 11).Value = 4000
$ws.Cells.Item(93, 13).Value = -1504

$ws.Cells.Item(96, 8).Value = 1312.5
$ws.Cells.Item(96, 9).Value = 1125
$ws.Cells.Item(96, 11).Value = 1125
$ws.Cells.Item(96, 13).Value = 248

$ws.Cells.Item(140, 8).Value = 99995
$ws.Cells.Item(140, 10).Value = 99995
$ws.Cells.Item(140, 12).Value = 99995
$ws.Cells.Item(140, 14).Value = -110355
